$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update stats for 2025-07 (row 20)
$ws.Range("B20").Value = 6165
$ws.Range("D20").Value = 5573138
$ws.Range("E20").Value = 903.9964314679643
$ws.Range("F20").Value = 6.49507686992572
$ws.Range("H20").Value = 26.06989931171113
